$d = $word.ActiveDocument

# --- Change 1: insert a new bullet paragraph before the "Pronađeni odgovarajući..." item ---
$count = $d.Paragraphs.Count
$targetIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Pronađeni odgovarajući*") {
        $targetIdx = $i
        break
    }
}
if ($targetIdx -eq -1) {
    throw "Could not find target paragraph for insertion"
}
$targetPara = $d.Paragraphs.Item($targetIdx)
$targetPara.Range.InsertParagraphBefore()

# The freshly inserted (empty) paragraph now sits at the same index the target
# paragraph used to occupy; fill it in with the new bullet text.
$newPara = $d.Paragraphs.Item($targetIdx)
$newPara.Range.Text = "Mala promjena slika – bila je greška zbog dimenzija"

# --- Change 2: split "Mjenjan hint ekran..." so "hint" becomes "help" in its own run ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("hint", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find 'hint' to replace"
}
$rng2.Text = "help"
# Toggle bold on/off to force Word to keep "help" as its own run, distinct from the
# surrounding "Mjenjan " / " ekran ..." text (mirrors how Word splits runs on a
# targeted in-place retype rather than a plain Find/Replace).
$rng2.Bold = 1
$rng2.Bold = 0

Write-Output "done"
